$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style (checkmark style s=4) from A19 into A20 before setting its value
$ws.Range("A19").Copy($ws.Range("A20"))

# Row 20
$ws.Range("A20").Value = '√'
$ws.Range("B20").Value = 'CQT2024083015-S002'
$ws.Range("C20").Value = '/cqproj22/Project/BD_CQT2024083015-F002_human_20240927'
$ws.Range("D20").Value = 45588
$ws.Range("E20").Value = '对F002的结果下，将cluster1，2，5抽出来， 进行亚群分析。核心目的还是为了找粒细胞。不需要做定义。输出下面基因气泡图就行:basophil：PTPRC/ITGAM/FCER1A/KIT/CD63/ENPP3;eosinophil：PTPRC/ITGAM/CCR3/SIGLEC8/FCGR3A/CLC/IL5RA/EPX'
$ws.Rows.Item(20).RowHeight = 14.25

# Row 21
$ws.Range("B21").Value = 'CQT2024040804-S005'
$ws.Range("C21").Value = '/cqproj22/Project/SC_CQT2024040804-F001_mouse_20240906'
$ws.Range("D21").Value = 45594
$ws.Range("E21").Value = '后面的分析都是针对S004的定义进行：
1.将小胶进行亚群分析，包括自动定义（CELLMARKER_auto_Microglia文件夹）和cluster层面的（直接以clusterN命名细胞类型，CELLMARKER_manual_Microglia文件夹），以及组间分析；
2. 若第1点的亚群分析获得的cluster数量较多(大于6个cluster)，请调整参数将cluster亚群数量控制在4～5个，再进行包括自动定义和cluster层面的，以及组间分析；同时第1点也要输出结果。如果默认参数获得的cluster数量在4～5个，则直接输出相应结果；
3. 所有组间分析的富集结果，均需要区分上下调基因和不区分上下调的基因的富集结果。提取组间比较的差异基因的富集结果中与能量代谢（或者是代谢）相关的通路(关键词“metabolism"搜索即可)，一个cluster的代谢通路画在同一张图上(bar图)，同时给出数据列表。
4. 做完前面亚群分析，在整合的小胶中，提取每个cluster与其他cluster比较的差异基因的富集结果中与能量代谢（或者是代谢）相关的通路(关键词“metabolism"搜索即可)，一个cluster的代谢通路画在同一张图上(bar图)，同时给出数据列表。
5. 针对S004的diff group文件夹的小胶差异基因做火山图，标记排名前十的基因，按照例图。
6. FABP5和AOPE在大群整合的、整合后分样本的表达热图；在小胶亚群分析中整合的、整合后分样本的表达热图。
7. 小胶亚群分析中：输出FABP5+小胶和AOPE+小胶在各样本的细胞数及占比。
（以上全部内容需要一起交付）'
$ws.Rows.Item(21).RowHeight = 14.25

# Row 22
$ws.Range("B22").Value = 'CQT2024041816'
$ws.Range("D22").Value = 45595
$ws.Range("E22").Value = '线上答疑'
$ws.Rows.Item(22).RowHeight = 14.25

# Row 23
$ws.Range("B23").Value = 'CQT2024062823-S002'
$ws.Range("C23").Value = '/cqproj22/Project/SC_CQT2024062823-F003_mouse_20240919'
$ws.Range("D23").Value = 45595
$ws.Range("E23").Value = '针对CQT2024062823-S001肌肉的结果：
需求1：对其中的0/1/3/9/10/13/14/19/22进行亚群分析，不需要定义，但是绘制总的marker基因气泡图
需求2：对其中的2/5/6/11/17/18/20进行亚群分析，不需要定义，但是绘制总的marker基因气泡图。
需求3：在大群上绘制一下基因的气泡图：
MTJ：Col1a1/Col3a1/Myod1/Myf5/Mmp2/Mmp9/Plin1/Fap
Adipocyte：Adipoq,Plin,Dgat2,Car3,Pparg
NMJ：Musk/Lrp/Vamp'
$ws.Rows.Item(23).RowHeight = 14.25

# Row 24
$ws.Range("B24").Value = 'CQT2024062823-F008'
$ws.Range("C24").Value = '/cqproj22/Project/SC_CQT2024062823-F008_mouse_20241029'
$ws.Range("D24").Value = 45596
$ws.Range("E24").Value = '6个肌肉样本的标准分析（数据来源：F001、F003、F007），客户有提供marker'
$ws.Rows.Item(24).RowHeight = 14.25

# Row 25
$ws.Range("B25").Value = 'CQT2023111808-S009'
$ws.Range("C25").Value = '/cqproj22/Project/SC_CQT2023111808-F002_mouse_20240709'
$ws.Range("D25").Value = 45597
$ws.Range("E25").Value = '对S008结果中中性粒的分群结果：计算KEGG中铁死亡（mmu04216）/凋亡（mmu04210）通路的score热图，umap展示（抽取通路中的基因，利用Seurat的AddModuleScore计算一个分数）；对S008结果中macro的分群结果：展示Mki67，Top2a在每个群的表达量umap热图和小提琴图，计算自噬通路（mmu04140）的score热图，Umap展示'
$ws.Rows.Item(25).RowHeight = 14.25

# Row 26
$ws.Range("B26").Value = 'CQT2024062823-S003'
$ws.Range("C26").Value = '/cqproj22/Project/SC_CQT2024062823-F003_mouse_20240919'
$ws.Range("D26").Value = 45597
$ws.Range("E26").Value = '基于CQT2024062823-F003结果，按照cluster进行分析： 绘制Sparcl1和Fndc5的小提琴图和umap热图。 需要每个cluster按照组别绘制和展示。 IM_1/2作为IM组，Young_1M/2M作为Young组'
$ws.Rows.Item(26).RowHeight = 14.25

# Row 27
$ws.Range("B27").Value = 'CQT2024101206-S001'
$ws.Range("C27").Value = '/cqproj22/Project/SC_CQT2024101206-S001_mouse_20241101'
$ws.Range("D27").Value = 45597
$ws.Range("E27").Value = '小鼠。需要查看是否捕获到病毒片段。病毒片段与病毒引物见附件。(构建加入病毒片段的新参考基因组/cqproj51/CellRanger_Reference/tmp_CQT2024101206-S001/mm10_wihtNewGenes)'
$ws.Rows.Item(27).RowHeight = 14.25

# Adjust column E width to match target (135.25 chars, closest achievable value)
$ws.Columns.Item(5).ColumnWidth = 134.53571428571428

# Update selection to match target state
$ws.Range("E27").Select()
